# Update the cryptos price list: refreshed Price (D) / Volume(1h) (E)
# figures for most rows, plus a 3-way reorder of the RenderToken /
# VeChain / InjectiveProtocol rows (49-51) with their own refreshed data.
#
# Cells whose new text happens to look like a number (e.g. "0.998",
# "7.00", "0.0240") are pre-formatted as Text ("@") so Excel keeps them
# as literal strings (preserving trailing zeros / precision) instead of
# silently coercing them to numeric values, matching the source data
# which stores every Price/Volume cell as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.720.47"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "2.599.98"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.07"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.32"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "2.622.13"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.46"
$ws.Range("E10").Value = "  -3.30%  "
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("E12").Value = "  -4.49%  "
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").Value = "3.062.45"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").Value = "60.478.04"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.25"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000140"
$ws.Range("E17").Value = "  +2.26%  "
$ws.Range("D18").Value = "2.612.44"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.39"
$ws.Range("E19").Value = "  +9.08%  "
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "346.26"
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.00"
$ws.Range("E22").Value = "  +8.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.535"
$ws.Range("E24").Value = "  +14.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.26"
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("E27").Value = "  -1.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.71"
$ws.Range("E28").Value = "  +3.88%  "
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("E30").Value = "  +9.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.48"
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("E35").Value = "  +3.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.974"
$ws.Range("E36").Value = "  +9.81%  "
$ws.Range("E37").Value = "  +3.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.82"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.82"
$ws.Range("E40").Value = "  +3.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.848"
$ws.Range("E41").Value = "  -3.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "295.59"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "138.02"
$ws.Range("E43").Value = "  +1.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.610"
$ws.Range("E45").Value = "  +1.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.71"
$ws.Range("E47").Value = "  +2.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0546"
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0240"
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.82"
$ws.Range("E50").Value = "  +4.99%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.90"
$ws.Range("E51").Value = "  +8.04%  "
